$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.651.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.133.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.21%  "
# Row 4
$ws.Range("E4").Value = "  -0.07%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.77%  "
# Row 7
$ws.Range("E7").Value = "  -0.16%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.86%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.145.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.117"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.25%  "
# Row 11
$ws.Range("E11").Value = "  -2.75%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.85%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.678.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.27%  "
# Row 14
$ws.Range("E14").Value = "  -2.18%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.640.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.14%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.138.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.32%  "
# Row 18
$ws.Range("E18").Value = "  -2.21%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "413.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.02%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.90%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.14%  "
# Row 22
$ws.Range("E22").Value = "  -1.57%  "
# Row 23
$ws.Range("E23").Value = "  -0.03%  "
# Row 24
$ws.Range("E24").Value = "  -2.30%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.481"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.56%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.194"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.62%  "
# Row 27
$ws.Range("E27").Value = "  -2.64%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.76%  "
# Row 29
$ws.Range("E29").Value = "  -0.28%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
# Row 31
$ws.Range("E31").Value = "  -1.52%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.54%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "163.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.61%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.92%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.53%  "
# Row 36
$ws.Range("E36").Value = "  -1.07%  "
# Row 37
$ws.Range("E37").Value = "  -1.38%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.26%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.617.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.53%  "
# Row 40
$ws.Range("E40").Value = "  -2.54%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.89%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "
# Row 43
$ws.Range("E43").Value = "  -3.71%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0613"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.77%  "
# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.47%  "
# Row 46
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "289.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0254"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.51%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.86%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
# Row 50
$ws.Range("E50").Value = "  -1.72%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.21%  "
